# Replace every "congenital" label in column A with "misc_long_term"
# across all affected worksheets (regression name-list workbook).

$wb = $excel.ActiveWorkbook

$targets = @(
    @{ Sheet = "variables_179"; Cell = "A2" },
    @{ Sheet = "variables_180"; Cell = "A2" },
    @{ Sheet = "variables_181"; Cell = "A2" },
    @{ Sheet = "variables_182"; Cell = "A2" },
    @{ Sheet = "variables_183"; Cell = "A2" },
    @{ Sheet = "variables_184"; Cell = "A2" },
    @{ Sheet = "variables_185"; Cell = "A2" },
    @{ Sheet = "variables_186"; Cell = "A2" },
    @{ Sheet = "variables_187"; Cell = "A2" },
    @{ Sheet = "variables_188"; Cell = "A2" },
    @{ Sheet = "variables_189"; Cell = "A2" },
    @{ Sheet = "variables_190"; Cell = "A2" },
    @{ Sheet = "variables_191"; Cell = "A2" },
    @{ Sheet = "variables_192"; Cell = "A2" },
    @{ Sheet = "variables_193"; Cell = "A2" },
    @{ Sheet = "variables_194"; Cell = "A2" },
    @{ Sheet = "variables_195"; Cell = "A2" },
    @{ Sheet = "variables_196"; Cell = "A2" },
    @{ Sheet = "variables_232"; Cell = "A5" },
    @{ Sheet = "variables_240"; Cell = "A4" },
    @{ Sheet = "variables_241"; Cell = "A4" },
    @{ Sheet = "variables_242"; Cell = "A4" },
    @{ Sheet = "variables_243"; Cell = "A4" },
    @{ Sheet = "variables_244"; Cell = "A4" },
    @{ Sheet = "variables_245"; Cell = "A4" },
    @{ Sheet = "variables_246"; Cell = "A4" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $ws.Range($t.Cell).Value = "misc_long_term"
}
